$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1=14, Q1=15 with the same (bold/bordered) style as the
# rest of row 1 (copy format only from O1, then set the values).
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: swap I/K and M/O column values, and add new P/Q columns.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I
    $ws.Cells.Item($r, 11).Value = 1  # K
    $ws.Cells.Item($r, 13).Value = 2  # M
    $ws.Cells.Item($r, 15).Value = 1  # O
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
